$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1031.7142
$ws.Range("I17").Value = 300
$ws.Range("J17").Value = 1153.6666
$ws.Range("K17").Value = 900
$ws.Range("L17").Value = 3460.9998
$ws.Range("M17").Value = -732
$ws.Range("N17").Value = -3796.9998
$ws.Range("H62").Value = 37043704
$ws.Range("I62").Value = 37043704
$ws.Range("K62").Value = 37043704
$ws.Range("M62").Value = -37043080
$ws.Range("H65").Value = 37043704
$ws.Range("I65").Value = 37043704
$ws.Range("K65").Value = 185218520
$ws.Range("M65").Value = -185215400
$ws.Range("H107").Value = 2847.2354
$ws.Range("I107").Value = 2699.889
$ws.Range("J107").Value = 3013
$ws.Range("K107").Value = 2699.889
$ws.Range("L107").Value = 3013
$ws.Range("M107").Value = -779.8890000000001
$ws.Range("N107").Value = -6853
$ws.Range("H129").Value = 810.54346
$ws.Range("J129").Value = 882.5789
$ws.Range("L129").Value = 2647.7367
$ws.Range("N129").Value = -12647.7367
$ws.Range("H137").Value = 1636.4828
$ws.Range("I137").Value = 1129.9
$ws.Range("J137").Value = 1903.1052
$ws.Range("K137").Value = 3389.7
$ws.Range("L137").Value = 5709.3156
$ws.Range("M137").Value = -839.7000000000003
$ws.Range("N137").Value = -10809.3156
$ws.Range("H138").Value = 573700.4399999999
$ws.Range("I138").Value = 1392.6842
$ws.Range("J138").Value = 764469.7
$ws.Range("K138").Value = 4178.0526
$ws.Range("L138").Value = 2293409.1
$ws.Range("M138").Value = 961.9474
$ws.Range("N138").Value = -2303689.1
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 144.2
$ws.Range("I5").Value = 105.25
$ws.Range("K5").Value = 105.25
$ws.Range("M5").Value = 6.75
$ws.Range("H26").Value = 611.75
$ws.Range("I26").Value = 611.75
$ws.Range("K26").Value = 611.75
$ws.Range("M26").Value = -281.75
$ws.Range("H36").Value = 3975.3333
$ws.Range("I36").Value = 3975.3333
$ws.Range("K36").Value = 3975.3333
$ws.Range("M36").Value = -3629.3333
$ws.Range("H45").Value = 1394.1538
$ws.Range("I45").Value = 1297.75
$ws.Range("J45").Value = 1548.4
$ws.Range("K45").Value = 1297.75
$ws.Range("L45").Value = 1548.4
$ws.Range("M45").Value = -920.75
$ws.Range("N45").Value = -2302.4
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 144.2
$ws.Range("I4").Value = 105.25
$ws.Range("K4").Value = 105.25
$ws.Range("M4").Value = 9.75
$ws.Range("H31").Value = 400
$ws.Range("I31").Value = 400
$ws.Range("K31").Value = 400
$ws.Range("M31").Value = -148
$ws.Range("H33").Value = 21
$ws.Range("I33").Value = 21
$ws.Range("K33").Value = 21
$ws.Range("M33").Value = 315
$ws.Range("H36").Value = 445.66666
$ws.Range("I36").Value = 445.66666
$ws.Range("K36").Value = 445.66666
$ws.Range("M36").Value = 88.33334000000002
$ws.Range("H37").Value = 550
$ws.Range("I37").Value = 550
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 550
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -413
$ws.Range("N37").ClearContents()
$ws.Range("H80").Value = 503
$ws.Range("J80").Value = 580.125
$ws.Range("L80").Value = 580.125
$ws.Range("N80").Value = -2576.125
$ws.Range("H83").Value = 503
$ws.Range("J83").Value = 580.125
$ws.Range("L83").Value = 2900.625
$ws.Range("N83").Value = -12884.625
$ws.Range("H99").Value = 100001160
$ws.Range("I99").Value = 142858240
$ws.Range("J99").Value = 1303.6666
$ws.Range("K99").Value = 142858240
$ws.Range("L99").Value = 1303.6666
$ws.Range("M99").Value = -142856742
$ws.Range("N99").Value = -4299.6666
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 295.8889
$ws.Range("J7").Value = 393.66666
$ws.Range("L7").Value = 393.66666
$ws.Range("N7").Value = -619.66666
$ws.Range("H31").Value = 1393.7667
$ws.Range("I31").Value = 1331.0769
$ws.Range("J31").Value = 1441.7059
$ws.Range("K31").Value = 1331.0769
$ws.Range("L31").Value = 1441.7059
$ws.Range("M31").Value = -1036.0769
$ws.Range("N31").Value = -2031.7059
$ws.Range("H34").Value = 1393.7667
$ws.Range("I34").Value = 1331.0769
$ws.Range("J34").Value = 1441.7059
$ws.Range("K34").Value = 1331.0769
$ws.Range("L34").Value = 1441.7059
$ws.Range("M34").Value = -1129.0769
$ws.Range("N34").Value = -1845.7059
$ws.Range("H135").Value = 34833.332
$ws.Range("J135").Value = 34833.332
$ws.Range("L135").Value = 34833.332
$ws.Range("N135").Value = -44973.332
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 551
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H34").Value = 10001907
$ws.Range("J34").Value = 12502322
$ws.Range("L34").Value = 37506966
$ws.Range("N34").Value = -37507134
$ws.Range("H39").Value = 2915.739
$ws.Range("J39").Value = 2831.524
$ws.Range("L39").Value = 8494.572
$ws.Range("N39").Value = -9082.572
$ws.Range("H131").Value = 17860414
$ws.Range("I131").Value = 90909440
$ws.Range("J131").Value = 3984.6
$ws.Range("K131").Value = 272728320
$ws.Range("L131").Value = 11953.8
$ws.Range("M131").Value = -272723280
$ws.Range("N131").Value = -22033.8
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7433.5
$ws.Range("I80").Value = 12500
$ws.Range("J80").Value = 4900.25
$ws.Range("K80").Value = 12500
$ws.Range("L80").Value = 4900.25
$ws.Range("M80").Value = -11502
$ws.Range("N80").Value = -6896.25
$ws.Range("H83").Value = 7433.5
$ws.Range("I83").Value = 12500
$ws.Range("J83").Value = 4900.25
$ws.Range("K83").Value = 62500
$ws.Range("L83").Value = 24501.25
$ws.Range("M83").Value = -57508
$ws.Range("N83").Value = -34485.25
$ws.Range("H93").Value = 30000
$ws.Range("J93").Value = 30000
$ws.Range("L93").Value = 30000
$ws.Range("N93").Value = -33744
$ws.Range("H113").Value = 1305.2142
$ws.Range("I113").Value = 1182.5
$ws.Range("J113").Value = 1468.8334
$ws.Range("K113").Value = 1182.5
$ws.Range("L113").Value = 1468.8334
$ws.Range("M113").Value = 987.5
$ws.Range("N113").Value = -5808.8334
$ws.Range("H132").Value = 2463.3438
$ws.Range("I132").Value = 2505.05
$ws.Range("J132").Value = 2393.8333
$ws.Range("K132").Value = 7515.150000000001
$ws.Range("L132").Value = 7181.499899999999
$ws.Range("M132").Value = -4985.150000000001
$ws.Range("N132").Value = -12241.4999
$ws.Range("H133").Value = 50307.8
$ws.Range("J133").Value = 50307.8
$ws.Range("L133").Value = 50307.8
$ws.Range("N133").Value = -60427.8
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 4396
$ws.Range("J31").Value = 4995
$ws.Range("L31").Value = 4995
$ws.Range("N31").Value = -5491
$ws.Range("H32").Value = 1075
$ws.Range("J32").Value = 1500
$ws.Range("L32").Value = 1500
$ws.Range("N32").Value = -2134
$ws.Range("H35").Value = 1575
$ws.Range("I35").Value = 433.33334
$ws.Range("J35").Value = 5000
$ws.Range("K35").Value = 433.33334
$ws.Range("L35").Value = 5000
$ws.Range("M35").Value = -97.33334000000002
$ws.Range("N35").Value = -5672
$ws.Range("H39").Value = 6059
$ws.Range("I39").Value = 6059
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 6059
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -5599
$ws.Range("N39").ClearContents()
$ws.Range("H41").Value = 5037
$ws.Range("J41").Value = 5037
$ws.Range("L41").Value = 5037
$ws.Range("N41").Value = -5913
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H48").Value = 3000
$ws.Range("I48").Value = 3000
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 3000
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = -2339
$ws.Range("N48").ClearContents()
$ws.Range("H70").Value = 39999
$ws.Range("J70").Value = 39999
$ws.Range("L70").Value = 39999
$ws.Range("N70").Value = -40539
$ws.Range("H73").Value = 39999
$ws.Range("J73").Value = 39999
$ws.Range("L73").Value = 39999
$ws.Range("N73").Value = -41871
$ws.Range("H136").Value = 2026
$ws.Range("I136").Value = 1521.6
$ws.Range("K136").Value = 4564.799999999999
$ws.Range("M136").Value = -2014.799999999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 27318
$ws.Range("I3").Value = 500
$ws.Range("J3").Value = 40727
$ws.Range("K3").Value = 500
$ws.Range("L3").Value = 40727
$ws.Range("M3").Value = -386
$ws.Range("N3").Value = -40955
$ws.Range("H63").Value = 15437.714
$ws.Range("J63").Value = 16973
$ws.Range("L63").Value = 16973
$ws.Range("N63").Value = -18221
$ws.Range("H66").Value = 15437.714
$ws.Range("J66").Value = 16973
$ws.Range("L66").Value = 50919
$ws.Range("N66").Value = -57159
$ws.Range("H107").Value = 484.2353
$ws.Range("I107").Value = 405.16666
$ws.Range("J107").Value = 674
$ws.Range("K107").Value = 1215.49998
$ws.Range("L107").Value = 2022
$ws.Range("M107").Value = 704.5000199999999
$ws.Range("N107").Value = -5862
$ws.Range("H132").Value = 2665.45
$ws.Range("I132").Value = 2805.8333
$ws.Range("J132").Value = 2244.3
$ws.Range("K132").Value = 8417.499899999999
$ws.Range("L132").Value = 6732.900000000001
$ws.Range("M132").Value = -5887.499899999999
$ws.Range("N132").Value = -11792.9
